$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 41
$ws.Range("C41").Value = 113
$ws.Range("E41").Value = 11056852

# Row 63
$ws.Range("C63").Value = 14362
$ws.Range("E63").Value = 36193311

# Row 64
$ws.Range("C64").Value = 5213
$ws.Range("E64").Value = 20442161

# Row 65
$ws.Range("C65").Value = 2017
$ws.Range("E65").Value = 13657300

# Row 70
$ws.Range("C70").Value = 15734
$ws.Range("E70").Value = 24684987

# Row 83
$ws.Range("C83").Value = 3415
$ws.Range("E83").Value = 115800208

# Row 91
$ws.Range("C91").Value = 151149
$ws.Range("E91").Value = 482480640

# Row 92
$ws.Range("C92").Value = 409212
$ws.Range("E92").Value = 1596566709

# Row 93
$ws.Range("C93").Value = 209624
$ws.Range("E93").Value = 1309639040

# Row 95
$ws.Range("C95").Value = 50792
$ws.Range("E95").Value = 933718057

# Row 96
$ws.Range("C96").Value = 17307
$ws.Range("E96").Value = 795758616

# Row 143
$ws.Range("C143").Value = 64958
$ws.Range("E143").Value = 373531948

# Row 146
$ws.Range("C146").Value = 4269
$ws.Range("E146").Value = 161502590

# Row 172
$ws.Range("C172").Value = 22702
$ws.Range("E172").Value = 44686673

# Row 184
$ws.Range("C184").Value = 68737
$ws.Range("E184").Value = 134191957
